# Apply commit "Add data for 2021-11-24" to the carjacking-by-neighborhood-by-month workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date.
$ws.Name = "Through 2021-11-16"

# Update the header/shared-string label for the current (partial) month column.
$ws.Range("B1").Value = "November 2021 (through November 16)"

# Updated / newly-added data cells (row = neighborhood, column = month).
$ws.Range("M3").Value   = 9
$ws.Range("AI3").Value  = 5
$ws.Range("AT3").Value  = 4
$ws.Range("M4").Value   = 5
$ws.Range("AT4").Value  = 9
$ws.Range("BE5").Value  = 5
$ws.Range("M8").Value   = 5
$ws.Range("B11").Value  = 2
$ws.Range("BE12").Value = 2
$ws.Range("B13").Value  = 4
$ws.Range("AT28").Value = 1
$ws.Range("AT32").Value = 2
$ws.Range("BE37").Value = 3
$ws.Range("B39").Value  = 2
$ws.Range("B45").Value  = 3
$ws.Range("BE48").Value = 3
$ws.Range("B64").Value  = 4
$ws.Range("M72").Value  = 2
$ws.Range("AI89").Value = 1
